# Trab. 2 de lab. pc2
# 1) "Data de entrega:" changes from 26/10 to 23/08.
# 2) Drop the stray _GoBack bookmark left around "Dicionários".

$d = $word.ActiveDocument

# --- 1. Update the delivery date -------------------------------------------------
# "26" and "/10" live in two separate runs in the original markup, but Find/Replace
# works against the document's flowed text, so it matches across the run boundary
# and leaves a single run behind with the new text, just like Word itself would.
$d.Content.Find.Execute("26/10", $true, $false, $false, $false, $false, $true, 1, `
                         $false, "23/08", 2) | Out-Null

# --- 2. Remove the leftover _GoBack bookmark --------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
